# Update employee mobile number (and name) for EMP005's row in Emp_Details.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to EMP005: Name -> "AYSHA", Mobile -> 6360581288
$ws.Range("B6").Value = "AYSHA"
$ws.Range("C6").Value = 6360581288

# Move the active selection to A6, matching the saved view state
[void]$ws.Range("A6").Select()
